$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '67.785.19'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.797.42'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.48%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '599.19'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '165.37'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +1.04%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.47'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("E12").Value = '  -1.72%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '35.85'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.44%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.433.68'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.53%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '3.797.64'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.54%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '67.815.25'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '18.40'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.48%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.114'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("E19").Value = '  +0.87%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '463.23'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("E21").Value = '  -2.09%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.700'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("E23").Value = '  -5.66%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '82.84'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("E25").Value = '  +0.42%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.10'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.56%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.04'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  -0.07%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '3.945.48'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("E30").Value = '  -1.93%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.47'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +3.25%  '
$ws.Range("E32").Value = '  -1.30%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '29.24'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("E34").Value = '  +0.01%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '9.04'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.43%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0995'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.988'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.77'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("E41").Value = '  +0.01%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '44.76'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.32%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '47.56'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("E45").Value = '  +0.48%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '150.66'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.63%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.38'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +9.24%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.36'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '27.48'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '400.02'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("E51").Value = '  +1.94%  '
